$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.686.47"
$ws.Range("E2").Value = "  -4.37%  "
$ws.Range("D3").Value = "3.670.74"
$ws.Range("E3").Value = "  -4.65%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'593.35"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'179.49"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("D7").Value = "3.662.67"
$ws.Range("E7").Value = "  -4.74%  "
$ws.Range("D8").Value = "'0.622"
$ws.Range("E8").Value = "  -6.77%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "'0.711"
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  -7.44%  "
$ws.Range("D12").Value = "'54.69"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("D13").Value = "'0.0000286"
$ws.Range("E13").Value = "  -10.75%  "
$ws.Range("D14").Value = "'10.31"
$ws.Range("E14").Value = "  -8.04%  "
$ws.Range("D15").Value = "4.270.91"
$ws.Range("E15").Value = "  -4.49%  "
$ws.Range("D16").Value = "3.686.70"
$ws.Range("E16").Value = "  -4.54%  "
$ws.Range("D17").Value = "'19.27"
$ws.Range("E17").Value = "  -6.35%  "
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("E19").Value = "  -6.84%  "
$ws.Range("D20").Value = "'12.68"
$ws.Range("E20").Value = "  -7.52%  "
$ws.Range("D21").Value = "67.691.02"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("D22").Value = "'406.16"
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("D23").Value = "'4.53"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "'87.73"
$ws.Range("E24").Value = "  -6.55%  "
$ws.Range("D25").Value = "'2.98"
$ws.Range("E25").Value = "  -8.34%  "
$ws.Range("D26").Value = "'12.70"
$ws.Range("E26").Value = "  -7.33%  "
$ws.Range("D27").Value = "'10.97"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'3.83"
$ws.Range("E28").Value = "  -6.32%  "
$ws.Range("D29").Value = "'6.05"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("D30").Value = "'9.37"
$ws.Range("E30").Value = "  -7.83%  "
$ws.Range("D31").Value = "'32.45"
$ws.Range("E31").Value = "  -7.17%  "
$ws.Range("D32").Value = "'7.27"
$ws.Range("E32").Value = "  -7.76%  "
$ws.Range("D33").Value = "'12.36"
$ws.Range("E33").Value = "  -8.05%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("D35").Value = "'43.22"
$ws.Range("E35").Value = "  -11.47%  "
$ws.Range("D36").Value = "'63.58"
$ws.Range("E36").Value = "  -8.09%  "
$ws.Range("D37").Value = "'589.46"
$ws.Range("E37").Value = "  -4.67%  "
$ws.Range("D38").Value = "0.0₃0868"
$ws.Range("E38").Value = "  -10.12%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'0.396"
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").Value = "'2.96"
$ws.Range("E44").Value = "  -8.80%  "
$ws.Range("D45").Value = "'0.0431"
$ws.Range("E45").Value = "  -7.28%  "
$ws.Range("D46").Value = "'2.85"
$ws.Range("E46").Value = "  -12.44%  "
$ws.Range("D47").Value = "'9.17"
$ws.Range("E47").Value = "  -8.86%  "
$ws.Range("D48").Value = "2.772.15"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").Value = "'0.133"
$ws.Range("E49").Value = "  -7.09%  "
$ws.Range("D50").Value = "'2.67"
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("D51").Value = "'3.12"
$ws.Range("E51").Value = "  -5.02%  "
